{"js": "// Update the worksheet date and each \"A\u00f7B=\" division expression in the\n// table to the new values described by the commit diff. Every source string\n// is unique in the document, so a straightforward search-and-replace per pair\n// is sufficient and keeps existing run formatting (font, size) intact.\nconst replacements = [\n  [\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"],\n  [\"905\u00f79=\", \"641\u00f73=\"],\n  [\"830\u00f79=\", \"130\u00f73=\"],\n  [\"319\u00f74=\", \"101\u00f73=\"],\n  [\"355\u00f75=\", \"185\u00f78=\"],\n  [\"894\u00f76=\", \"522\u00f75=\"],\n  [\"505\u00f75=\", \"946\u00f78=\"],\n  [\"862\u00f75=\", \"927\u00f73=\"],\n  [\"660\u00f73=\", \"350\u00f79=\"],\n  [\"428\u00f79=\", \"700\u00f75=\"],\n  [\"360\u00f73=\", \"939\u00f72=\"],\n  [\"699\u00f78=\", \"660\u00f78=\"],\n  [\"633\u00f72=\", \"154\u00f74=\"],\n  [\"491\u00f78=\", \"639\u00f78=\"],\n  [\"815\u00f75=\", \"670\u00f79=\"],\n  [\"993\u00f77=\", \"778\u00f76=\"],\n  [\"245\u00f72=\", \"420\u00f77=\"],\n  [\"814\u00f79=\", \"807\u00f74=\"],\n  [\"406\u00f75=\", \"128\u00f73=\"],\n  [\"642\u00f76=\", \"375\u00f79=\"],\n  [\"697\u00f74=\", \"137\u00f78=\"],\n  [\"787\u00f75=\", \"209\u00f76=\"],\n  [\"819\u00f74=\", \"389\u00f76=\"],\n  [\"897\u00f73=\", \"283\u00f74=\"],\n  [\"894\u00f73=\", \"778\u00f75=\"],\n  [\"254\u00f74=\", \"652\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace every occurrence found (should be exactly one, since all source\n  // strings are unique within the document).\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each \"A\u00f7B=\" division expression in the\n# table to the new values described by the commit diff. Every source string\n# is unique in the document, so a straightforward Find/Replace per pair is\n# sufficient and preserves the existing run formatting (font, size).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-10-19 Saturday', '2024-10-20 Sunday'),\n    @('905\u00f79=', '641\u00f73='),\n    @('830\u00f79=', '130\u00f73='),\n    @('319\u00f74=', '101\u00f73='),\n    @('355\u00f75=', '185\u00f78='),\n    @('894\u00f76=', '522\u00f75='),\n    @('505\u00f75=', '946\u00f78='),\n    @('862\u00f75=', '927\u00f73='),\n    @('660\u00f73=', '350\u00f79='),\n    @('428\u00f79=', '700\u00f75='),\n    @('360\u00f73=', '939\u00f72='),\n    @('699\u00f78=', '660\u00f78='),\n    @('633\u00f72=', '154\u00f74='),\n    @('491\u00f78=', '639\u00f78='),\n    @('815\u00f75=', '670\u00f79='),\n    @('993\u00f77=', '778\u00f76='),\n    @('245\u00f72=', '420\u00f77='),\n    @('814\u00f79=', '807\u00f74='),\n    @('406\u00f75=', '128\u00f73='),\n    @('642\u00f76=', '375\u00f79='),\n    @('697\u00f74=', '137\u00f78='),\n    @('787\u00f75=', '209\u00f76='),\n    @('819\u00f74=', '389\u00f76='),\n    @('897\u00f73=', '283\u00f74='),\n    @('894\u00f73=', '778\u00f75='),\n    @('254\u00f74=', '652\u00f78='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
